$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Tracker sheet: append today's (Oct 06 2022 / serial 44840) closing row (32)
# ---------------------------------------------------------------------------
$tracker = $wb.Worksheets.Item("Tracker")

# NOTE: a whole-row Range.Copy / PasteSpecial(xlPasteFormats) from row 31
# onto row 32 here confuses this engine's recalculation of the
# COUNTA(Tracker!...)-driven formulas on other sheets (it leaves them
# showing stale pre-edit values even after the values below are written),
# so the day's numbers/formulas are entered directly and formatting is
# copied field-by-field afterwards instead.
$tracker.Range("A32").Value2 = 44840
$tracker.Range("B32").Value2 = 3771.97
$tracker.Range("C32").Value2 = 3744.52
$tracker.Range("D32").Value2 = 3769.4580542802919
$tracker.Range("E32").Value2 = 0.26200000000000001
$tracker.Range("F32").Value2 = 0
$tracker.Range("G32").Value2 = 10.08
$tracker.Range("H32").Value2 = 10.08
$tracker.Range("I32").Formula = "=G32/H32"
$tracker.Range("J32").Formula = "=J31+G32*100"
$tracker.Range("K32").Formula = "=H32*100+K31"

# Match the look of the row above (mirrors a user copying the row down
# before editing it) without using Copy/PasteSpecial.
$tracker.Range("A32").NumberFormat = $tracker.Range("A31").NumberFormat
$tracker.Range("B32:D32").NumberFormat = $tracker.Range("B31").NumberFormat
$tracker.Range("E32").NumberFormat = $tracker.Range("E31").NumberFormat
$tracker.Range("F32").NumberFormat = $tracker.Range("F31").NumberFormat
$tracker.Range("G32").NumberFormat = $tracker.Range("G31").NumberFormat
$tracker.Range("H32").NumberFormat = $tracker.Range("H31").NumberFormat
$tracker.Range("I32").NumberFormat = $tracker.Range("I31").NumberFormat
$tracker.Range("J32").NumberFormat = $tracker.Range("J31").NumberFormat
$tracker.Range("K32").NumberFormat = $tracker.Range("K31").NumberFormat

# ---------------------------------------------------------------------------
# Indicator sheet: today's manually-entered inputs (open/close quote, qty)
# ---------------------------------------------------------------------------
$indicator = $wb.Worksheets.Item("Indicator")
$indicator.Range("C2").Value2 = 1625
$indicator.Range("C8").Value2 = 3771.97
$indicator.Range("C9").Value2 = 3744.52

# ---------------------------------------------------------------------------
# Restore the end-of-session selections on each sheet
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$tracker.Activate()
$tracker.Range("J36").Select()

$indicator.Activate()
$indicator.Range("E18").Select()

$dashboard.Activate()
$dashboard.Range("M21").Select()

$indicator.Activate()
